$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1274.2142
$ws.Range("J41").Value = 1549.5454
$ws.Range("L41").Value = 1549.5454
$ws.Range("N41").Value = -2429.5454
$ws.Range("H70").Value = 6669.364
$ws.Range("J70").Value = 10227.167
$ws.Range("L70").Value = 30681.501
$ws.Range("N70").Value = -31221.501
$ws.Range("H73").Value = 6669.364
$ws.Range("J73").Value = 10227.167
$ws.Range("L73").Value = 30681.501
$ws.Range("N73").Value = -32553.501
$ws.Range("H80").Value = 11558.5
$ws.Range("I80").Value = 666.3333
$ws.Range("J80").Value = 16226.571
$ws.Range("K80").Value = 1998.9999
$ws.Range("L80").Value = 48679.713
$ws.Range("M80").Value = -1000.9999
$ws.Range("N80").Value = -50675.713
$ws.Range("H83").Value = 11558.5
$ws.Range("I83").Value = 666.3333
$ws.Range("J83").Value = 16226.571
$ws.Range("K83").Value = 5996.9997
$ws.Range("L83").Value = 146039.139
$ws.Range("M83").Value = -1004.9997
$ws.Range("N83").Value = -156023.139
$ws.Range("H86").Value = 2425.6667
$ws.Range("I86").Value = 2070.75
$ws.Range("K86").Value = 2070.75
$ws.Range("M86").Value = -947.75
$ws.Range("H89").Value = 2425.6667
$ws.Range("I89").Value = 2070.75
$ws.Range("K89").Value = 10353.75
$ws.Range("M89").Value = -4737.75
$ws.Range("H112").Value = 13527.6
$ws.Range("J112").Value = 19043.715
$ws.Range("L112").Value = 57131.145
$ws.Range("N112").Value = -59347.145
$ws.Range("H116").Value = 15485.694
$ws.Range("I116").Value = 16787
$ws.Range("K116").Value = 16787
$ws.Range("M116").Value = -13345
$ws.Range("H132").Value = 17455.182
$ws.Range("I132").Value = 17455.182
$ws.Range("K132").Value = 52365.546
$ws.Range("M132").Value = -49835.546
$ws.Range("H135").Value = 4513.9287
$ws.Range("I135").Value = 3899.375
$ws.Range("K135").Value = 35094.375
$ws.Range("M135").Value = -32559.375
$ws.Range("H138").Value = 40819.42
$ws.Range("I138").Value = 2005.7059
$ws.Range("J138").Value = 114134.22
$ws.Range("K138").Value = 6017.1177
$ws.Range("L138").Value = 342402.66
$ws.Range("M138").Value = -877.1176999999998
$ws.Range("N138").Value = -352682.66

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24212.043
$ws.Range("I32").Value = 27043.707
$ws.Range("J32").Value = 992.4
$ws.Range("K32").Value = 27043.707
$ws.Range("L32").Value = 992.4
$ws.Range("M32").Value = -26756.707
$ws.Range("N32").Value = -1566.4
$ws.Range("H45").Value = 2524.2307
$ws.Range("I45").Value = 1444.0588
$ws.Range("K45").Value = 1444.0588
$ws.Range("M45").Value = -1067.0588
$ws.Range("H61").Value = 8755
$ws.Range("I61").Value = 1250.5454
$ws.Range("K61").Value = 1250.5454
$ws.Range("M61").Value = -1038.5454
$ws.Range("H74").Value = 236066.81
$ws.Range("I74").Value = 316860.78
$ws.Range("J74").Value = 16768.857
$ws.Range("K74").Value = 316860.78
$ws.Range("L74").Value = 16768.857
$ws.Range("M74").Value = -315986.78
$ws.Range("N74").Value = -18516.857
$ws.Range("H77").Value = 236066.81
$ws.Range("I77").Value = 316860.78
$ws.Range("J77").Value = 16768.857
$ws.Range("K77").Value = 1584303.9
$ws.Range("L77").Value = 83844.285
$ws.Range("M77").Value = -1579935.9
$ws.Range("N77").Value = -92580.285
$ws.Range("H132").Value = 1474.8422
$ws.Range("I132").Value = 1279
$ws.Range("K132").Value = 3837
$ws.Range("M132").Value = -1307
$ws.Range("H136").Value = 8755
$ws.Range("I136").Value = 1250.5454
$ws.Range("K136").Value = 3751.6362
$ws.Range("M136").Value = -1201.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 12306.556
$ws.Range("J64").Value = 15612.571
$ws.Range("L64").Value = 15612.571
$ws.Range("N64").Value = -16062.571
$ws.Range("H67").Value = 12306.556
$ws.Range("J67").Value = 15612.571
$ws.Range("L67").Value = 15612.571
$ws.Range("N67").Value = -17172.571
$ws.Range("H80").Value = 889.53845
$ws.Range("I80").Value = 1147.8572
$ws.Range("J80").Value = 588.1667
$ws.Range("K80").Value = 1147.8572
$ws.Range("L80").Value = 588.1667
$ws.Range("M80").Value = -149.8571999999999
$ws.Range("N80").Value = -2584.1667
$ws.Range("H83").Value = 889.53845
$ws.Range("I83").Value = 1147.8572
$ws.Range("J83").Value = 588.1667
$ws.Range("K83").Value = 5739.286
$ws.Range("L83").Value = 2940.8335
$ws.Range("M83").Value = -747.2860000000001
$ws.Range("N83").Value = -12924.8335
$ws.Range("H105").Value = 1512.7667
$ws.Range("J105").Value = 1868.8572
$ws.Range("L105").Value = 1868.8572
$ws.Range("N105").Value = -5362.8572
$ws.Range("H132").Value = 86593
$ws.Range("J132").Value = 86593
$ws.Range("L132").Value = 86593
$ws.Range("N132").Value = -96713
$ws.Range("H134").Value = 12257.577
$ws.Range("I134").Value = 19361.143
$ws.Range("K134").Value = 58083.429
$ws.Range("M134").Value = -55548.429

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6668280
$ws.Range("I31").Value = 7143727.5
$ws.Range("K31").Value = 7143727.5
$ws.Range("M31").Value = -7143432.5
$ws.Range("H34").Value = 6668280
$ws.Range("I34").Value = 7143727.5
$ws.Range("K34").Value = 7143727.5
$ws.Range("M34").Value = -7143525.5
$ws.Range("H62").Value = 5175.2383
$ws.Range("I62").Value = 5247.778
$ws.Range("J62").Value = 5120.8335
$ws.Range("K62").Value = 5247.778
$ws.Range("L62").Value = 5120.8335
$ws.Range("M62").Value = -4623.778
$ws.Range("N62").Value = -6368.8335
$ws.Range("H65").Value = 5175.2383
$ws.Range("I65").Value = 5247.778
$ws.Range("J65").Value = 5120.8335
$ws.Range("K65").Value = 26238.89
$ws.Range("L65").Value = 25604.1675
$ws.Range("M65").Value = -23118.89
$ws.Range("N65").Value = -31844.1675
$ws.Range("H107").Value = 748.2381
$ws.Range("I107").Value = 692.46155
$ws.Range("J107").Value = 838.875
$ws.Range("K107").Value = 692.46155
$ws.Range("L107").Value = 838.875
$ws.Range("M107").Value = 1227.53845
$ws.Range("N107").Value = -4678.875
$ws.Range("H122").Value = 2130.6428
$ws.Range("I122").Value = 2022.7
$ws.Range("K122").Value = 6068.1
$ws.Range("M122").Value = -3618.1
$ws.Range("H134").Value = 2744.353
$ws.Range("I134").Value = 2603.4666
$ws.Range("K134").Value = 7810.399800000001
$ws.Range("M134").Value = -5275.399800000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1676.3
$ws.Range("J132").Value = 2406.8333
$ws.Range("L132").Value = 21661.4997
$ws.Range("N132").Value = -26721.4997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2590.65
$ws.Range("I102").Value = 2322.9443
$ws.Range("K102").Value = 2322.9443
$ws.Range("M102").Value = -700.9443000000001
$ws.Range("H126").Value = 1582.0667
$ws.Range("I126").Value = 1616.5
$ws.Range("J126").Value = 1100
$ws.Range("K126").Value = 4849.5
$ws.Range("L126").Value = 3300
$ws.Range("M126").Value = -2379.5
$ws.Range("N126").Value = -8240
$ws.Range("H132").Value = 2111.5
$ws.Range("I132").Value = 2090.8572
$ws.Range("J132").Value = 2256
$ws.Range("K132").Value = 6272.571599999999
$ws.Range("L132").Value = 6768
$ws.Range("M132").Value = -3742.571599999999
$ws.Range("N132").Value = -11828

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1678.6471
$ws.Range("I40").Value = 1623.875
$ws.Range("K40").Value = 1623.875
$ws.Range("M40").Value = -1487.875
$ws.Range("H93").Value = 1496.7222
$ws.Range("I93").Value = 1153.75
$ws.Range("J93").Value = 2182.6667
$ws.Range("K93").Value = 1153.75
$ws.Range("L93").Value = 2182.6667
$ws.Range("M93").Value = 94.25
$ws.Range("N93").Value = -4678.6667
$ws.Range("H105").Value = 65306.5
$ws.Range("J105").Value = 65306.5
$ws.Range("L105").Value = 65306.5
$ws.Range("N105").Value = -72294.5
$ws.Range("H132").Value = 6110.6665
$ws.Range("I132").Value = 6199.2
$ws.Range("K132").Value = 18597.6
$ws.Range("M132").Value = -16067.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 49591.668
$ws.Range("J51").Value = 95000
$ws.Range("L51").Value = 95000
$ws.Range("N51").Value = -96020
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H81").Value = 26332.666
$ws.Range("I81").Value = 26332.666
$ws.Range("K81").Value = 52665.332
$ws.Range("M81").Value = -51604.332
$ws.Range("H84").Value = 26332.666
$ws.Range("I84").Value = 26332.666
$ws.Range("K84").Value = 263326.66
$ws.Range("M84").Value = -258022.66
$ws.Range("H92").Value = 35000
$ws.Range("J92").Value = 35000
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -39992
$ws.Range("H95").Value = 84865
$ws.Range("J95").Value = 84865
$ws.Range("L95").Value = 84865
$ws.Range("N95").Value = -90357
$ws.Range("H126").Value = 3482.1482
$ws.Range("I126").Value = 2913.348
$ws.Range("K126").Value = 8740.044
$ws.Range("M126").Value = -6270.044
$ws.Range("H132").Value = 91925
$ws.Range("I132").Value = 143480
$ws.Range("K132").Value = 430440
$ws.Range("M132").Value = -427910
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
